$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking cells in column D so they are
# preserved as text (matching the workbooks existing text-encoded prices)
# rather than being auto-converted to numbers by Excel.

$ws.Range("D2").Value = "27.148.62"
$ws.Range("E2").Value = "  -2.53%  "

$ws.Range("D3").Value = "1.868.66"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.53"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5136"
$ws.Range("E7").Value = "  +2.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3748"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8846"
$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.66"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "1.876.84"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07585"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.335"
$ws.Range("E14").Value = "  -2.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.15"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008552"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.14"
$ws.Range("E18").Value = "  -2.63%  "

$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "27.192.69"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.039"
$ws.Range("E21").Value = "  -2.57%  "

$ws.Range("D22").Value = "2.117.15"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  -1.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.473"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.54"
$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.842"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.03"
$ws.Range("E27").Value = "  -1.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.171"
$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.01"
$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.736"
$ws.Range("E30").Value = "  -3.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.697"
$ws.Range("E31").Value = "  +1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09041"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05156"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.098"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7569"
$ws.Range("E35").Value = "  -0.93%  "

$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02039"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.522"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.029"
$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5387"
$ws.Range("E41").Value = "  -3.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.660"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.23"
$ws.Range("E43").Value = "  +3.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.543"
$ws.Range("E44").Value = "  +0.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1486"
$ws.Range("E45").Value = "  -1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4681"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.12%  "

$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.575"
$ws.Range("E49").Value = "  -3.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.99"
$ws.Range("E50").Value = "  -3.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.46"
$ws.Range("E51").Value = "  -1.54%  "
